$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts "Klant" and everything below it down by one)
$ws.Rows.Item(5).Insert()

# The engine's row-insert doesn't always carry the column-A border format onto
# the freshly inserted row, so restore it explicitly (thin right border,
# matching the style used by every other data row in column A).
$ws.Range("A5").Borders.Item(10).Color = 0
$ws.Range("A5").Borders.Item(10).Weight = 2
$ws.Range("A5").Borders.Item(10).LineStyle = 1

# Populate the new row with the "2e projectleider" column mapping
$ws.Range("A5").Value = "2e projectleider"
$ws.Range("B5").Value = "Output"
$ws.Range("C5").Value = "2e Projectleider"
$ws.Range("G5").Value = "2e Projectleider"

# Update the active selection to match the saved workbook state
[void]$ws.Range("F11").Select()
